$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "newSheet"

# Update numeric values in columns C-F for rows 2-16
$ws.Range("C2").Value = 4399.13
$ws.Range("D2").Value = 133.31
$ws.Range("E2").Value = 4998.47
$ws.Range("F2").Value = 151.47
$ws.Range("C3").Value = 4792.82
$ws.Range("D3").Value = 145.24
$ws.Range("E3").Value = 564.22
$ws.Range("F3").Value = 17.1
$ws.Range("C4").Value = 393.69
$ws.Range("D4").Value = 11.93
$ws.Range("E4").Value = -4434.25
$ws.Range("F4").Value = -134.37
$ws.Range("C5").Value = 4485.64
$ws.Range("D5").Value = 135.93
$ws.Range("E5").Value = 5084.7
$ws.Range("F5").Value = 154.08
$ws.Range("C6").Value = 982.0700000000001
$ws.Range("D6").Value = 29.76
$ws.Range("E6").Value = 1147.3
$ws.Range("F6").Value = 34.77
$ws.Range("C7").Value = -3503.57
$ws.Range("D7").Value = -106.17
$ws.Range("E7").Value = -3937.4
$ws.Range("F7").Value = -119.31
$ws.Range("C8").Value = 4871.58
$ws.Range("D8").Value = 147.62
$ws.Range("E8").Value = 4649.2
$ws.Range("F8").Value = 140.88
$ws.Range("C9").Value = 2880.57
$ws.Range("D9").Value = 87.29000000000001
$ws.Range("E9").Value = 2859.85
$ws.Range("F9").Value = 86.66
$ws.Range("C10").Value = -1991.01
$ws.Range("D10").Value = -60.33
$ws.Range("E10").Value = -1789.35
$ws.Range("F10").Value = -54.22
$ws.Range("C11").Value = 1266.52
$ws.Range("D11").Value = 38.38
$ws.Range("E11").Value = 2759.79
$ws.Range("F11").Value = 83.63
$ws.Range("C12").Value = 2789.06
$ws.Range("D12").Value = 84.52
$ws.Range("E12").Value = 1857.1
$ws.Range("F12").Value = 56.28
$ws.Range("C13").Value = 1522.54
$ws.Range("D13").Value = 46.14
$ws.Range("E13").Value = -902.6900000000001
$ws.Range("F13").Value = -27.35
$ws.Range("C14").Value = 744.79
$ws.Range("D14").Value = 22.57
$ws.Range("E14").Value = 798.35
$ws.Range("F14").Value = 24.19
$ws.Range("C15").Value = 794.78
$ws.Range("D15").Value = 24.08
$ws.Range("E15").Value = 869.01
$ws.Range("F15").Value = 26.33
$ws.Range("C16").Value = 49.99
$ws.Range("D16").Value = 1.51
$ws.Range("E16").Value = 70.66
$ws.Range("F16").Value = 2.14
